# Generate Report for Handback
#
# For each of the two "handback" detail sheets (zh-cn, de-de):
#   - Column F ("Latest Target File")   gets the same value/link as column A (source file)
#   - Column G ("Latest Handback File") gets the same value/link as column D (latest handoff file)
#   - Column H ("Latest Handback DateTime") gets stamped with the handback time
# The shared "Status" text moves from "Ready for handoff" to
# "Handed back: in sync with en-US" for every row that uses it.

function Find-HyperlinkAddress($ws, $addrTarget) {
    foreach ($hl in @($ws.Hyperlinks)) {
        $addr = $hl.Range.Address()
        if ($addr -eq $addrTarget) {
            return $hl.Address
        }
    }
    return $null
}

$wb = $excel.ActiveWorkbook

# Status text moves from "Ready for handoff" to "Handed back: in sync with
# en-US" for every data row on both locale sheets.
foreach ($sheetName in @("zh-cn", "de-de")) {
    $statusWs = $wb.Worksheets.Item($sheetName)
    foreach ($row in @(2, 3)) {
        $statusWs.Range("C" + $row).Value = "Handed back: in sync with en-US"
    }
}

$handbackTimes = @{
    "zh-cn" = @{ "2" = "2016-03-21 00:40:01"; "3" = "2016-03-21 00:40:01" }
    "de-de" = @{ "2" = "2016-03-21 00:40:09"; "3" = "2016-03-21 00:40:09" }
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in @(2, 3)) {
        $aAddr = "`$A`$" + $row
        $dAddr = "`$D`$" + $row
        $fCell = "F" + $row
        $gCell = "G" + $row
        $hCell = "H" + $row

        $aValue = $ws.Range("A" + $row).Value()
        $dValue = $ws.Range("D" + $row).Value()
        $aLink = Find-HyperlinkAddress $ws $aAddr
        $dLink = Find-HyperlinkAddress $ws $dAddr

        # F: Latest Target File - mirrors the source file, now that it is in sync
        $ws.Hyperlinks.Add($ws.Range($fCell), $aLink, "", "", $aValue)
        $ws.Range($fCell).Style = $ws.Range("A" + $row).Style

        # G: Latest Handback File - mirrors the file handed off, now handed back
        $ws.Hyperlinks.Add($ws.Range($gCell), $dLink, "", "", $dValue)
        $ws.Range($gCell).Style = $ws.Range("D" + $row).Style

        # H: Latest Handback DateTime
        $ws.Range($hCell).Value = $handbackTimes[$sheetName][[string]$row]
    }
}
